$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the date in A1 (45308 -> 45309, i.e. 1/17/2024 -> 1/18/2024) ---
$ws.Range("A1").Value = 45309

# --- Update the three prices in column D ---
$ws.Range("D33").Value = 1165.8
$ws.Range("D34").Value = 1966.2
$ws.Range("D35").Value = 2164.56
